$wb = $excel.ActiveWorkbook

# --- Sheet "Castle" (sheet1) ---------------------------------------------
$castle = $wb.Worksheets.Item("Castle")

# Clear the old placeholder grid (A1:C3) so stale cells don't linger.
$castle.Range("A1:C3").ClearContents() | Out-Null

# Header row.
$castle.Cells.Item(1,1).Value = "ID"
$castle.Cells.Item(1,2).Value = "Castle Name"
$castle.Cells.Item(1,3).Value = "Faction"
$castle.Cells.Item(1,4).Value = "Wall Strength"
$castle.Cells.Item(1,5).Value = "Troops"

# Data row.
$castle.Cells.Item(2,2).Value = "Castle at Old Town"
$castle.Cells.Item(2,3).Value = "Knights of the Round"
$castle.Cells.Item(2,4).Value = 20000
$castle.Cells.Item(2,5).Value = "Mark"
$castle.Cells.Item(2,6).Value = "Talison"
$castle.Cells.Item(2,7).Value = "/"

# Column widths (engine quantizes ColumnWidth to 1/6-character increments,
# so these are the closest achievable settings to the authored widths).
$castle.Columns.Item(2).ColumnWidth = 20.666666666666668
$castle.Columns.Item(3).ColumnWidth = 25.666666666666668
$castle.Columns.Item(4).ColumnWidth = 12.333333333333334

$castle.Range("F3").Select() | Out-Null

# --- Sheet "Troop" (sheet2) -----------------------------------------------
$troop = $wb.Worksheets.Item("Troop")

# Clear the old placeholder grid (A1:C3) so stale cells don't linger.
$troop.Range("A1:C3").ClearContents() | Out-Null

# Header row.
$troop.Cells.Item(1,1).Value = "ID"
$troop.Cells.Item(1,2).Value = "Troop Name"
$troop.Cells.Item(1,3).Value = "Salute"
$troop.Cells.Item(1,4).Value = "Strength"

# Data rows.
$troop.Cells.Item(2,2).Value = "Mark"
$troop.Cells.Item(2,3).Value = "Yes sir?"
$troop.Cells.Item(2,4).Value = 12

$troop.Cells.Item(3,2).Value = "Talison"
$troop.Cells.Item(3,3).Value = "Yes sir!"
$troop.Cells.Item(3,4).Value = "asdf"

# Column widths.
$troop.Columns.Item(1).ColumnWidth = 4.5
$troop.Columns.Item(2).ColumnWidth = 11.0
$troop.Columns.Item(3).ColumnWidth = 14.5
$troop.Columns.Item(4).ColumnWidth = 7.666666666666667

$troop.Range("D3").Select() | Out-Null

# Troop becomes the active/visible tab (activeTab=1, tabSelected moves here).
$troop.Activate()
